# Applies the "CrimeController.php class" update described by the diff:
#  - fills in the empty "Classe: " / variable-list paragraphs for the
#    Controller package with CrimeController + natureza/tempo/ano vars
#  - moves the _GoBack bookmark (Word re-stamps it at the last edit
#    location) from the old empty paragraph to the new trailing one
#  - shifts the lastRenderedPageBreak marker from "crimeVW" onto the
#    "Classe: persistence.php" run above it (consequence of the new
#    content pushing the page break earlier)
#
# Note: Paragraph.Range.Text always includes the trailing paragraph
# mark, so text comparisons below use -like "...*" rather than -eq.

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the (now stale) lastRenderedPageBreak before "crimeVW".
# ---------------------------------------------------------------------
$pCrimeVW = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "crimeVW*") { $pCrimeVW = $p }
}
$xmlCrimeVW = "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='2'/></w:numPr><w:rPr><w:b/><w:sz w:val='32'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:proofErr w:type='gramStart'/><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t>crimeVW</w:t></w:r><w:proofErr w:type='spellEnd'/><w:proofErr w:type='gramEnd'/><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t xml:space='preserve'> -&gt; </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t>crimeView</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$pCrimeVW.Range.InsertXML($xmlCrimeVW)

# ---------------------------------------------------------------------
# 2) Add lastRenderedPageBreak before "Classe: persistence.php".
# ---------------------------------------------------------------------
$pPersistence = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Classe: persistence.php*") { $pPersistence = $p }
}
$xmlPersistence = "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:rPr><w:sz w:val='28'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Classe: </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t>persistence.php</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$pPersistence.Range.InsertXML($xmlPersistence)

# ---------------------------------------------------------------------
# 3) Drop the _GoBack bookmark from its old resting place (the empty,
#    underlined paragraph right before "Pacote: Test").
# ---------------------------------------------------------------------
$pOldBookmark = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.WordOpenXML -like "*_GoBack*") { $pOldBookmark = $p }
}
$xmlOldBookmark = "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:rPr><w:b/><w:sz w:val='28'/><w:u w:val='single'/></w:rPr></w:pPr></w:p>"
$pOldBookmark.Range.InsertXML($xmlOldBookmark)

# ---------------------------------------------------------------------
# 4) Fill in "Classe: " -> "Classe: CrimeController" (Controller package).
# ---------------------------------------------------------------------
$pClasse = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Classe: *" -and $p.Range.Text.Trim() -eq "Classe:") { $pClasse = $p }
}
$xmlClasse = "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:sz w:val='28'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t xml:space='preserve'>Classe: </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='28'/></w:rPr><w:t>CrimeController</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$pClasse.Range.InsertXML($xmlClasse)

# ---------------------------------------------------------------------
# 5) Expand the following empty list paragraph into the three new
#    member-variable mappings, plus a trailing empty paragraph that now
#    carries the _GoBack bookmark.
# ---------------------------------------------------------------------
$pVars = $pClasse.Next()
$xmlVars  = "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='gramStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>natureza</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t xml:space='preserve'> -&gt; </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>kindCrime</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$xmlVars += "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='gramStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>tempo</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t xml:space='preserve'> -&gt; </w:t></w:r><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>time</w:t></w:r></w:p>"
$xmlVars += "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='gramStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>ano</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t xml:space='preserve'> -&gt; </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>year</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$xmlVars += "<w:p $w><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:ind w:left='1440'/><w:rPr><w:sz w:val='28'/></w:rPr></w:pPr><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$pVars.Range.InsertXML($xmlVars)
